$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The D/E (and, for the two swapped rows, B/C) columns hold
# text-formatted values (prices like "38.061.12" and percentages like
# "  +2.90%  "), not real numbers/percents. Force each touched cell to
# Text format right before writing its new value so Excel does not
# silently coerce the assigned strings into numeric values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.065.64'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.052.55'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.59'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.614'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.93'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +7.76%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.386'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0813'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.22%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.355.92'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.62'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.94'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.752'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.48%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.053.01'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.956.60'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.76%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.76'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0836'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.51'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.07%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.22'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.27'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.28%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.02'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.32'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.21%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.54'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.07%  '
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.58'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.23%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.06'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +11.02%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.08'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +13.71%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.36%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.536.85'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.77'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.37%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.89'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.74'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.62%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.88%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.04'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +13.50%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.44%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.12'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.242.09'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.37%  '
